$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5348869.5
$ws.Range("J17").Value = 5348869.5
$ws.Range("L17").Value = 16046608.5
$ws.Range("N17").Value = -16046944.5
$ws.Range("H40").Value = 2010.581
$ws.Range("I40").Value = 2007.0807
$ws.Range("K40").Value = 2007.0807
$ws.Range("M40").Value = -1832.0807
$ws.Range("H74").Value = 5128.364
$ws.Range("I74").Value = 4333.3335
$ws.Range("J74").Value = 6082.4
$ws.Range("K74").Value = 4333.3335
$ws.Range("L74").Value = 6082.4
$ws.Range("M74").Value = -3397.3335
$ws.Range("N74").Value = -7954.4
$ws.Range("H77").Value = 5128.364
$ws.Range("I77").Value = 4333.3335
$ws.Range("J77").Value = 6082.4
$ws.Range("K77").Value = 21666.6675
$ws.Range("L77").Value = 30412
$ws.Range("M77").Value = -16986.6675
$ws.Range("N77").Value = -39772
$ws.Range("H112").Value = 13290129
$ws.Range("I112").Value = 333.33334
$ws.Range("J112").Value = 14286864
$ws.Range("K112").Value = 1000.00002
$ws.Range("L112").Value = 42860592
$ws.Range("M112").Value = 107.9999799999999
$ws.Range("N112").Value = -42862808
$ws.Range("H129").Value = 1133.9272
$ws.Range("I129").Value = 760
$ws.Range("J129").Value = 1171.32
$ws.Range("K129").Value = 2280
$ws.Range("L129").Value = 3513.96
$ws.Range("M129").Value = 2720
$ws.Range("N129").Value = -13513.96

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 254.44444
$ws.Range("I4").Value = 181.66667
$ws.Range("K4").Value = 181.66667
$ws.Range("M4").Value = -65.66667000000001
$ws.Range("H32").Value = 5830.35
$ws.Range("I32").Value = 4952.053
$ws.Range("J32").Value = 19590.334
$ws.Range("K32").Value = 4952.053
$ws.Range("L32").Value = 19590.334
$ws.Range("M32").Value = -4665.053
$ws.Range("N32").Value = -20164.334
$ws.Range("H61").Value = 308777.38
$ws.Range("I61").Value = 6293.8076
$ws.Range("J61").Value = 1432287.8
$ws.Range("K61").Value = 6293.8076
$ws.Range("L61").Value = 1432287.8
$ws.Range("M61").Value = -6081.8076
$ws.Range("N61").Value = -1432711.8
$ws.Range("H74").Value = 1960.3871
$ws.Range("I74").Value = 1535.7778
$ws.Range("J74").Value = 2548.3076
$ws.Range("K74").Value = 1535.7778
$ws.Range("L74").Value = 2548.3076
$ws.Range("M74").Value = -661.7778000000001
$ws.Range("N74").Value = -4296.3076
$ws.Range("H77").Value = 1960.3871
$ws.Range("I77").Value = 1535.7778
$ws.Range("J77").Value = 2548.3076
$ws.Range("K77").Value = 7678.889
$ws.Range("L77").Value = 12741.538
$ws.Range("M77").Value = -3310.889
$ws.Range("N77").Value = -21477.538
$ws.Range("H110").Value = 3433.0833
$ws.Range("I110").Value = 3562.182
$ws.Range("J110").Value = 2013
$ws.Range("K110").Value = 3562.182
$ws.Range("L110").Value = 2013
$ws.Range("M110").Value = -1517.182
$ws.Range("N110").Value = -6103
$ws.Range("H132").Value = 2567361
$ws.Range("I132").Value = 2363.68
$ws.Range("J132").Value = 7147713
$ws.Range("K132").Value = 7091.039999999999
$ws.Range("L132").Value = 21443139
$ws.Range("M132").Value = -4561.039999999999
$ws.Range("N132").Value = -21448199
$ws.Range("H136").Value = 308777.38
$ws.Range("I136").Value = 6293.8076
$ws.Range("J136").Value = 1432287.8
$ws.Range("K136").Value = 18881.4228
$ws.Range("L136").Value = 4296863.4
$ws.Range("M136").Value = -16331.4228
$ws.Range("N136").Value = -4301963.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H99").Value = 55557420
$ws.Range("I99").Value = 90910200
$ws.Range("J99").Value = 3054.2856
$ws.Range("K99").Value = 90910200
$ws.Range("L99").Value = 3054.2856
$ws.Range("M99").Value = -90908702
$ws.Range("N99").Value = -6050.2856
$ws.Range("H107").Value = 1970.3334
$ws.Range("I107").Value = 2255.5
$ws.Range("J107").Value = 1400
$ws.Range("K107").Value = 2255.5
$ws.Range("L107").Value = 1400
$ws.Range("M107").Value = -335.5
$ws.Range("N107").Value = -5240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4338.615
$ws.Range("I31").Value = 2010.4651
$ws.Range("J31").Value = 8889.091
$ws.Range("K31").Value = 2010.4651
$ws.Range("L31").Value = 8889.091
$ws.Range("M31").Value = -1715.4651
$ws.Range("N31").Value = -9479.091
$ws.Range("H34").Value = 4338.615
$ws.Range("I34").Value = 2010.4651
$ws.Range("J34").Value = 8889.091
$ws.Range("K34").Value = 2010.4651
$ws.Range("L34").Value = 8889.091
$ws.Range("M34").Value = -1808.4651
$ws.Range("N34").Value = -9293.091
$ws.Range("H134").Value = 199784.8
$ws.Range("I134").Value = 3967.394
$ws.Range("J134").Value = 558783.4
$ws.Range("K134").Value = 11902.182
$ws.Range("L134").Value = 1676350.2
$ws.Range("M134").Value = -9367.181999999999
$ws.Range("N134").Value = -1681420.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1556.8422
$ws.Range("I34").Value = 850
$ws.Range("J34").Value = 1745.3334
$ws.Range("K34").Value = 2550
$ws.Range("L34").Value = 5236.0002
$ws.Range("M34").Value = -2466
$ws.Range("N34").Value = -5404.0002
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H103").Value = 6187
$ws.Range("I103").Value = 7424.2856
$ws.Range("J103").Value = 3300
$ws.Range("K103").Value = 22272.8568
$ws.Range("L103").Value = 9900
$ws.Range("M103").Value = -21393.8568
$ws.Range("N103").Value = -11658
$ws.Range("H117").Value = 23819274
$ws.Range("I117").Value = 33709.668
$ws.Range("J117").Value = 30306246
$ws.Range("K117").Value = 101129.004
$ws.Range("L117").Value = 90918738
$ws.Range("M117").Value = -97687.00399999999
$ws.Range("N117").Value = -90925622
$ws.Range("H129").Value = 1238.7333
$ws.Range("I129").Value = 591.6667
$ws.Range("J129").Value = 1670.1111
$ws.Range("K129").Value = 1775.0001
$ws.Range("L129").Value = 5010.3333
$ws.Range("M129").Value = 3224.9999
$ws.Range("N129").Value = -15010.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5344.3887
$ws.Range("I70").Value = 5346.6
$ws.Range("J70").Value = 5333.3335
$ws.Range("K70").Value = 5346.6
$ws.Range("L70").Value = 5333.3335
$ws.Range("M70").Value = -5076.6
$ws.Range("N70").Value = -5873.3335
$ws.Range("H73").Value = 5344.3887
$ws.Range("I73").Value = 5346.6
$ws.Range("J73").Value = 5333.3335
$ws.Range("K73").Value = 5346.6
$ws.Range("L73").Value = 5333.3335
$ws.Range("M73").Value = -4410.6
$ws.Range("N73").Value = -7205.3335
$ws.Range("H123").Value = 14772.444
$ws.Range("J123").Value = 14772.444
$ws.Range("L123").Value = 14772.444
$ws.Range("N123").Value = -19672.444
$ws.Range("H132").Value = 6437.515
$ws.Range("I132").Value = 7397.391
$ws.Range("J132").Value = 4229.8
$ws.Range("K132").Value = 22192.173
$ws.Range("L132").Value = 12689.4
$ws.Range("M132").Value = -19662.173
$ws.Range("N132").Value = -17749.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 55558052
$ws.Range("I40").Value = 66669184
$ws.Range("J40").Value = 2400
$ws.Range("K40").Value = 66669184
$ws.Range("L40").Value = 2400
$ws.Range("M40").Value = -66669048
$ws.Range("N40").Value = -2672
$ws.Range("H61").Value = 1800
$ws.Range("I61").Value = 1800
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1800
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1598
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 1800
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 370
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 5502028.5
$ws.Range("I122").Value = 5960322.5
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 17880967.5
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -17878517.5
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2426.7856
$ws.Range("I132").Value = 1880.125
$ws.Range("K132").Value = 5640.375
$ws.Range("M132").Value = -3110.375
$ws.Range("H136").Value = 2668.92
$ws.Range("I136").Value = 2959.24
$ws.Range("J136").Value = 2378.6
$ws.Range("K136").Value = 8877.719999999999
$ws.Range("L136").Value = 7135.799999999999
$ws.Range("M136").Value = -6327.719999999999
$ws.Range("N136").Value = -12235.8
